# Updated remaining queries for C3DC
# - Fix JOIN conditions in all SQL queries (std.id/prt.id -> std.study_id/prt.participant_id,
#   with matching qualified column renames on the right-hand side of each ON clause)
# - Move the active selection to B2
# - Widen column C (drop "best fit" auto-width, set an explicit wider width)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells on Sheet1 that contain the embedded SQL text needing the JOIN-column fix.
$sqlCells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $sqlCells) {
    $rng = $ws.Range($addr)
    $text = [string]$rng.Value2

    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

    $rng.Value = $text
}

# Move the active cell selection from C7 to B2.
$ws.Range("B2").Select() | Out-Null

# Column C: remove the "best fit" auto width and set an explicit, slightly wider width.
$ws.Columns.Item(3).ColumnWidth = 66.5
